$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.416292
$ws.Range("H2").Value = 43.248876
$ws.Range("I2").Value = 0.8004770782290026
$ws.Range("J2").Value = 0.8004770782290026
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 68.00339
$ws.Range("N2").Value = 204.01017
$ws.Range("O2").Value = 0.6265962299909886
$ws.Range("P2").Value = 0.6265962299909885
$ws.Range("Q2").Value = 980.3567272298799
$ws.Range("R2").Value = 8823.210545068918
$ws.Range("S2").Value = 0.5015759194124947
$ws.Range("T2").Value = 0.5015759194124946

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.416292
$ws.Range("H3").Value = 43.248876
$ws.Range("I3").Value = 0.8004770782290026
$ws.Range("J3").Value = 0.8004770782290026
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.991529999999999
$ws.Range("N3").Value = 26.97459
$ws.Range("O3").Value = 0.08284967558015671
$ws.Range("P3").Value = 0.08284967558015671
$ws.Range("Q3").Value = 129.62452200676
$ws.Range("R3").Value = 1166.62069806084
$ws.Range("S3").Value = 0.06631926624062459
$ws.Range("T3").Value = 0.06631926624062459

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.416292
$ws.Range("H4").Value = 43.248876
$ws.Range("I4").Value = 0.8004770782290026
$ws.Range("J4").Value = 0.8004770782290026
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.76843933333333
$ws.Range("N4").Value = 32.305318
$ws.Range("O4").Value = 0.09922245772090688
$ws.Range("P4").Value = 0.09922245772090688
$ws.Range("Q4").Value = 155.2409658136187
$ws.Range("R4").Value = 1397.168692322568
$ws.Range("S4").Value = 0.07942530305113228
$ws.Range("T4").Value = 0.07942530305113228

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 14.416292
$ws.Range("H5").Value = 43.248876
$ws.Range("I5").Value = 0.8004770782290026
$ws.Range("J5").Value = 0.8004770782290026
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 20.764887
$ws.Range("N5").Value = 62.294661
$ws.Range("O5").Value = 0.1913316367079478
$ws.Range("P5").Value = 0.1913316367079478
$ws.Range("Q5").Value = 299.352674339004
$ws.Range("R5").Value = 2694.174069051036
$ws.Range("S5").Value = 0.153156589524751
$ws.Range("T5").Value = 0.153156589524751

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.178646666666667
$ws.Range("H6").Value = 9.53594
$ws.Range("I6").Value = 0.1764971045575167
$ws.Range("J6").Value = 0.1764971045575167
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 68.00339
$ws.Range("N6").Value = 204.01017
$ws.Range("O6").Value = 0.6265962299909886
$ws.Range("P6").Value = 0.6265962299909885
$ws.Range("Q6").Value = 216.1587489455333
$ws.Range("R6").Value = 1945.4287405098
$ws.Range("S6").Value = 0.1105924203200653
$ws.Range("T6").Value = 0.1105924203200653

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.178646666666667
$ws.Range("H7").Value = 9.53594
$ws.Range("I7").Value = 0.1764971045575167
$ws.Range("J7").Value = 0.1764971045575167
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.991529999999999
$ws.Range("N7").Value = 26.97459
$ws.Range("O7").Value = 0.08284967558015671
$ws.Range("P7").Value = 0.08284967558015671
$ws.Range("Q7").Value = 28.58089686273333
$ws.Range("R7").Value = 257.2280717646
$ws.Range("S7").Value = 0.01462272785342726
$ws.Range("T7").Value = 0.01462272785342726

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.178646666666667
$ws.Range("H8").Value = 9.53594
$ws.Range("I8").Value = 0.1764971045575167
$ws.Range("J8").Value = 0.1764971045575167
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.76843933333333
$ws.Range("N8").Value = 32.305318
$ws.Range("O8").Value = 0.09922245772090688
$ws.Range("P8").Value = 0.09922245772090688
$ws.Range("Q8").Value = 34.22906379210222
$ws.Range("R8").Value = 308.06157412892
$ws.Range("S8").Value = 0.01751247649482069
$ws.Range("T8").Value = 0.01751247649482068

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.178646666666667
$ws.Range("H9").Value = 9.53594
$ws.Range("I9").Value = 0.1764971045575167
$ws.Range("J9").Value = 0.1764971045575167
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.764887
$ws.Range("N9").Value = 62.294661
$ws.Range("O9").Value = 0.1913316367079478
$ws.Range("P9").Value = 0.1913316367079478
$ws.Range("Q9").Value = 66.00423884625999
$ws.Range("R9").Value = 594.03814961634
$ws.Range("S9").Value = 0.03376947988920347
$ws.Range("T9").Value = 0.03376947988920347

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Mfng"
$ws.Range("C10").Value = "Notch1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.4146863333333333
$ws.Range("H10").Value = 1.244059
$ws.Range("I10").Value = 0.02302581721348076
$ws.Range("J10").Value = 0.02302581721348076
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 68.00339
$ws.Range("N10").Value = 204.01017
$ws.Range("O10").Value = 0.6265962299909886
$ws.Range("P10").Value = 0.6265962299909885
$ws.Range("Q10").Value = 28.20007645333666
$ws.Range("R10").Value = 253.80068808003
$ws.Range("S10").Value = 0.01442789025842865
$ws.Range("T10").Value = 0.01442789025842865

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Mfng"
$ws.Range("C11").Value = "Notch1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.4146863333333333
$ws.Range("H11").Value = 1.244059
$ws.Range("I11").Value = 0.02302581721348076
$ws.Range("J11").Value = 0.02302581721348076
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.991529999999999
$ws.Range("N11").Value = 26.97459
$ws.Range("O11").Value = 0.08284967558015671
$ws.Range("P11").Value = 0.08284967558015671
$ws.Range("Q11").Value = 3.728664606756666
$ws.Range("R11").Value = 33.55798146081
$ws.Range("S11").Value = 0.001907681486104869
$ws.Range("T11").Value = 0.001907681486104869

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Mfng"
$ws.Range("C12").Value = "Notch1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.4146863333333333
$ws.Range("H12").Value = 1.244059
$ws.Range("I12").Value = 0.02302581721348076
$ws.Range("J12").Value = 0.02302581721348076
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 10.76843933333333
$ws.Range("N12").Value = 32.305318
$ws.Range("O12").Value = 0.09922245772090688
$ws.Range("P12").Value = 0.09922245772090688
$ws.Range("Q12").Value = 4.465524622862445
$ws.Range("R12").Value = 40.189721605762
$ws.Range("S12").Value = 0.002284678174953924
$ws.Range("T12").Value = 0.002284678174953924

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Mfng"
$ws.Range("C13").Value = "Notch1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.4146863333333333
$ws.Range("H13").Value = 1.244059
$ws.Range("I13").Value = 0.02302581721348076
$ws.Range("J13").Value = 0.02302581721348076
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 20.764887
$ws.Range("N13").Value = 62.294661
$ws.Range("O13").Value = 0.1913316367079478
$ws.Range("P13").Value = 0.1913316367079478
$ws.Range("Q13").Value = 8.610914852111
$ws.Range("R13").Value = 77.49823366899899
$ws.Range("S13").Value = 0.004405567293993311
$ws.Range("T13").Value = 0.004405567293993311
